$wb = $excel.ActiveWorkbook
$tx = $wb.Worksheets.Item("Transaksi")
$sm = $wb.Worksheets.Item("Summary")

# ---------------------------------------------------------------------------
# Helper: write a "formatted number" text cell (stored as text, e.g. "0,00")
# ---------------------------------------------------------------------------
function Set-Text($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

# Helper: write a date-serial cell, keeping the yyyy-mm-dd date format
function Set-DateCell($ws, $addr, $serial) {
    $ws.Range($addr).Value = $serial
    $ws.Range($addr).NumberFormat = "YYYY-MM-DD"
}

# Helper: blank out a cell entirely (no value, no style) like the target's
# bare <c t="inlineStr"/> cells
function Set-Blank($ws, $addr) {
    $ws.Range($addr).ClearContents()
    $ws.Range($addr).ClearFormats()
}

# ---------------------------------------------------------------------------
# Sheet "Transaksi"
# ---------------------------------------------------------------------------

# Row 2 (Opening Balance) - date moved out of A2, balances bumped
Set-Blank    $tx "A2"
Set-Text     $tx "D2" "496.264.489,00"
Set-DateCell $tx "E2" 45657
Set-Text     $tx "H2" "496.264.489,00"

# Row 3
Set-DateCell $tx "A3" 45658
Set-Text     $tx "B3" "0,00"
Set-Text     $tx "C3" "10.000,00"
Set-Text     $tx "D3" "496.254.489,00"
Set-DateCell $tx "E3" 45658
Set-Text     $tx "F3" "10.000,00"
Set-Text     $tx "G3" "0,00"
Set-Text     $tx "H3" "496.254.489,00"

# Row 4
Set-DateCell $tx "A4" 45659
Set-Text     $tx "C4" "20.138,00"
Set-Text     $tx "D4" "496.234.351,00"
Set-DateCell $tx "E4" 45659
Set-Text     $tx "F4" "20.138,00"
Set-Text     $tx "H4" "496.234.351,00"

# Row 5
Set-DateCell $tx "A5" 45660
Set-Text     $tx "B5" "157.751.823,00"
Set-Text     $tx "C5" "0,00"
Set-Text     $tx "D5" "653.986.174,00"
Set-DateCell $tx "E5" 45660
Set-Text     $tx "F5" "0,00"
Set-Text     $tx "G5" "157.751.823,00"
Set-Text     $tx "H5" "653.986.174,00"

# Row 6
Set-DateCell $tx "A6" 45667
Set-Text     $tx "C6" "150.025.000,00"
Set-Text     $tx "D6" "503.961.174,00"
Set-Blank    $tx "E6"
Set-Text     $tx "F6" "0,00"
Set-Text     $tx "H6" "653.986.174,00"
Set-Text     $tx "J6" "150.025.000,00"
Set-Text     $tx "K6" "Unmatched"
Set-Text     $tx "L6" "G1"

# Row 7
Set-Blank    $tx "A7"
Set-Text     $tx "B7" "0,00"
Set-Text     $tx "D7" "503.961.174,00"
Set-DateCell $tx "E7" 45667
Set-Text     $tx "F7" "150.000.000,00"
Set-Text     $tx "G7" "0,00"
Set-Text     $tx "H7" "503.986.174,00"
Set-Text     $tx "J7" "-150.000.000,00"
Set-Text     $tx "K7" "Unmatched"
Set-Text     $tx "L7" "G1"

# Row 8
Set-Blank    $tx "A8"
Set-Text     $tx "B8" "0,00"
Set-Text     $tx "C8" "0,00"
Set-Text     $tx "D8" "503.961.174,00"
Set-DateCell $tx "E8" 45667
Set-Text     $tx "F8" "25.000,00"
Set-Text     $tx "G8" "0,00"
Set-Text     $tx "H8" "503.961.174,00"
Set-Text     $tx "J8" "-25.000,00"
Set-Text     $tx "K8" "Unmatched"
Set-Text     $tx "L8" "G1"

# Row 9 (new)
Set-DateCell $tx "A9" 45688
Set-Text     $tx "B9" "431.953,00"
Set-Text     $tx "C9" "0,00"
Set-Text     $tx "D9" "504.393.127,00"
Set-DateCell $tx "E9" 45688
Set-Text     $tx "F9" "0,00"
Set-Text     $tx "G9" "431.953,00"
Set-Text     $tx "H9" "504.393.127,00"
Set-Text     $tx "I9" "0,00"
Set-Text     $tx "J9" "0,00"
Set-Text     $tx "K9" "Matched"
Set-Text     $tx "L9" "-"

# Row 10 (new)
Set-DateCell $tx "A10" 45688
Set-Text     $tx "B10" "0,00"
Set-Text     $tx "C10" "86.391,00"
Set-Text     $tx "D10" "504.306.736,00"
Set-DateCell $tx "E10" 45688
Set-Text     $tx "F10" "86.391,00"
Set-Text     $tx "G10" "0,00"
Set-Text     $tx "H10" "504.306.736,00"
Set-Text     $tx "I10" "0,00"
Set-Text     $tx "J10" "0,00"
Set-Text     $tx "K10" "Matched"
Set-Text     $tx "L10" "-"

# Row 11 (new)
Set-DateCell $tx "A11" 45688
Set-Text     $tx "B11" "0,00"
Set-Text     $tx "C11" "12.000,00"
Set-Text     $tx "D11" "504.294.736,00"
Set-DateCell $tx "E11" 45688
Set-Text     $tx "F11" "12.000,00"
Set-Text     $tx "G11" "0,00"
Set-Text     $tx "H11" "504.294.736,00"
Set-Text     $tx "I11" "0,00"
Set-Text     $tx "J11" "0,00"
Set-Text     $tx "K11" "Matched"
Set-Text     $tx "L11" "-"

# Row 12 (new)
Set-DateCell $tx "A12" 45688
Set-Text     $tx "B12" "0,00"
Set-Text     $tx "C12" "10.000,00"
Set-Text     $tx "D12" "504.284.736,00"
Set-DateCell $tx "E12" 45688
Set-Text     $tx "F12" "10.000,00"
Set-Text     $tx "G12" "0,00"
Set-Text     $tx "H12" "504.284.736,00"
Set-Text     $tx "I12" "0,00"
Set-Text     $tx "J12" "0,00"
Set-Text     $tx "K12" "Matched"
Set-Text     $tx "L12" "-"

# Row 13 (new, Closing Balance)
Set-Text     $tx "B13" "158.183.776,00"
Set-Text     $tx "C13" "150.163.529,00"
Set-Text     $tx "D13" "504.284.736,00"
Set-DateCell $tx "E13" 45688
Set-Text     $tx "F13" "150.163.529,00"
Set-Text     $tx "G13" "158.183.776,00"
Set-Text     $tx "H13" "504.284.736,00"
Set-Text     $tx "I13" "0,00"
Set-Text     $tx "J13" "0,00"
Set-Text     $tx "K13" "Closing Balance"

# ---------------------------------------------------------------------------
# Sheet "Summary"
# ---------------------------------------------------------------------------
Set-Blank $sm "B2"
Set-Text  $sm "C2" "496.264.489,00"
Set-Text  $sm "D2" "496.264.489,00"

Set-Blank $sm "B3"
Set-Text  $sm "C3" "504.284.736,00"
Set-Text  $sm "D3" "504.284.736,00"
